# ==========================================================================
# Edit script: rename "Metadata" -> "MetaData" (new sheetId), rewrite the
# ATM sheet summary rows and the MetaData sheet content, and restyle the
# blank separator rows on the Calls sheet.
# ==========================================================================

$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------
# 1) Calls sheet: blank separator rows get the "Normal" style re-applied
#    (cosmetic re-stamp matching the authoring tool's re-save; harmless
#    no-op visually, default formatting either way).
# --------------------------------------------------------------------
$callsWs = $wb.Worksheets.Item("Calls")
$blankRows = 27,33,39,45,51,57,63,69,75,81
foreach ($r in $blankRows) {
    $rng = $callsWs.Range("A" + $r + ":G" + $r)
    $rng.Font.Bold = $false
}

# --------------------------------------------------------------------
# 2) Price sheet: just a view/selection change.
# --------------------------------------------------------------------
$priceWs = $wb.Worksheets.Item("Price")
$priceWs.Range("A24").Select()

# --------------------------------------------------------------------
# 3) ATM sheet: replace the single option-chain row with the three
#    summary lines describing the ATM call identified for the trade.
# --------------------------------------------------------------------
$atmWs = $wb.Worksheets.Item("ATM")
$atmWs.Cells.Clear()

$atmWs.Range("A1").Value = "CAT 10/17/25 C470"
$atmWs.Range("A1").Font.Bold = $true
$atmWs.Range("A1").IndentLevel = 1
$atmWs.Range("A1").HorizontalAlignment = -4131
$atmWs.Range("A1").VerticalAlignment = -4108

$prefix2 = "Trade-date underlying = "
$bold2 = "466.54"
$suffix2 = " on 19 Sep 2025."
$atmWs.Range("A2").Value = $prefix2 + $bold2 + $suffix2
$start2 = $prefix2.Length + 1
$len2 = $bold2.Length
$atmWs.Range("A2").Characters($start2, $len2).Font.Bold = $true

$atmWs.Range("A3").Value = "Daily closing prices of CAT 10/17/25 C470 from 19 Sep 2025 to 17 Oct 2025 inclusive."
$atmWs.Range("A3").Font.Bold = $true
$atmWs.Range("A3").IndentLevel = 1
$atmWs.Range("A3").HorizontalAlignment = -4131
$atmWs.Range("A3").VerticalAlignment = -4108

$atmWs.Range("A9").Select()

# --------------------------------------------------------------------
# 4) Metadata sheet: delete the old free-form notes sheet and add a
#    freshly-created "MetaData" sheet (so it gets a brand-new sheetId)
#    in the same tab position, holding the cleaned-up 4-point summary.
# --------------------------------------------------------------------
$atmAnchor = $wb.Worksheets.Item("ATM")
$newMeta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $atmAnchor)
$newMeta.Name = "MetaData_TEMP_NAME"

$oldMeta = $wb.Worksheets.Item("Metadata")
$oldMeta.Delete()

$metaWs = $wb.Worksheets.Item("MetaData_TEMP_NAME")
$metaWs.Name = "MetaData"

$q1 = [char]0x201C + "Closing prices for all options (all of the strikes available for puts and calls) on your assigned trade date of 19th September 2025 with an expiry date of 17th October 2025." + [char]0x201D
$q2 = [char]0x201C + "Closing prices for the At The Money (ATM) call option (strike as close as possible to closing price of stock) on your trade date." + [char]0x201D
$q3 = [char]0x201C + "Prices for that same option strike for each date from 19th September 2025 to 17th October 2025 inclusive (required to delta hedge the option)." + [char]0x201D
$q4 = [char]0x201C + "Closing prices for your stock for matching dates and a history of prices leading up to the trade date to estimate historical volatility." + [char]0x201D

# --- Point 1 ---
$metaWs.Range("A1").Value = "1. Full option chain on the trade date"
$metaWs.Range("A1").Font.Bold = $true

$metaWs.Range("A3").Value = $q1
$metaWs.Range("A3").IndentLevel = 1
$metaWs.Range("A3").HorizontalAlignment = -4131
$metaWs.Range("A3").VerticalAlignment = -4108

# --- Point 2 ---
$metaWs.Range("A5").Value = "2.  ATM call price on the trade date"
$metaWs.Range("A5").Font.Bold = $true

$metaWs.Range("A7").Value = $q2
$metaWs.Range("A7").IndentLevel = 1
$metaWs.Range("A7").HorizontalAlignment = -4131
$metaWs.Range("A7").VerticalAlignment = -4108

# --- Point 3 (rich text: "same" in bold italic) ---
$p3a = "3. Daily prices for that "
$p3b = "same"
$p3c = " ATM call from 19 Sep to 17 Oct"
$metaWs.Range("A9").Value = $p3a + $p3b + $p3c
$metaWs.Range("A9").Font.Bold = $true
$b3start = $p3a.Length + 1
$b3len = $p3b.Length
$c3start = $b3start + $b3len
$c3len = $p3c.Length
$metaWs.Range("A9").Characters($b3start, $b3len).Font.Italic = $true
$metaWs.Range("A9").Characters($c3start, $c3len).Font.Bold = $true

$metaWs.Range("A11").Value = $q3
$metaWs.Range("A11").IndentLevel = 1
$metaWs.Range("A11").HorizontalAlignment = -4131
$metaWs.Range("A11").VerticalAlignment = -4108

# --- Point 4 ---
$metaWs.Range("A13").Value = "4. Stock prices for matching dates + history before trade date"
$metaWs.Range("A13").Font.Bold = $true

$metaWs.Range("A15").Value = $q4
$metaWs.Range("A15").IndentLevel = 1
$metaWs.Range("A15").HorizontalAlignment = -4131
$metaWs.Range("A15").VerticalAlignment = -4108

$metaWs.Range("I28").Select()

$wb.Worksheets.Item("MetaData").Activate()
